$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '69.386.53'
Set-TextValue $ws.Range("E2") '  -1.12%  '

Set-TextValue $ws.Range("D3") '3.537.75'
Set-TextValue $ws.Range("E3") '  -1.82%  '

Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.27%  '

Set-TextValue $ws.Range("D5") '195.94'
Set-TextValue $ws.Range("E5") '  -0.41%  '

Set-TextValue $ws.Range("D6") '583.17'
Set-TextValue $ws.Range("E6") '  -3.55%  '

Set-TextValue $ws.Range("D7") '0.609'
Set-TextValue $ws.Range("E7") '  -2.68%  '

Set-TextValue $ws.Range("D8") '1.00'
Set-TextValue $ws.Range("E8") '  +0.02%  '

Set-TextValue $ws.Range("D9") '0.203'
Set-TextValue $ws.Range("E9") '  -1.80%  '

Set-TextValue $ws.Range("D10") '0.629'
Set-TextValue $ws.Range("E10") '  -2.95%  '

Set-TextValue $ws.Range("D11") '51.69'
Set-TextValue $ws.Range("E11") '  -3.96%  '

Set-TextValue $ws.Range("E12") '  -6.30%  '

Set-TextValue $ws.Range("D13") '9.21'
Set-TextValue $ws.Range("E13") '  -3.67%  '

Set-TextValue $ws.Range("D14") '4.104.80'
Set-TextValue $ws.Range("E14") '  -1.73%  '

Set-TextValue $ws.Range("D15") '665.33'
Set-TextValue $ws.Range("E15") '  +12.25%  '

Set-TextValue $ws.Range("D16") '69.545.78'
Set-TextValue $ws.Range("E16") '  -1.07%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D17") '12.51'
Set-TextValue $ws.Range("E17") '  -4.53%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D18") '3.537.99'
Set-TextValue $ws.Range("E18") '  -1.78%  '

Set-TextValue $ws.Range("E19") '  -0.80%  '

Set-TextValue $ws.Range("D20") '18.45'
Set-TextValue $ws.Range("E20") '  -3.62%  '

Set-TextValue $ws.Range("D21") '0.963'
Set-TextValue $ws.Range("E21") '  -3.26%  '

Set-TextValue $ws.Range("D22") '18.35'
Set-TextValue $ws.Range("E22") '  +3.41%  '

Set-TextValue $ws.Range("D23") '5.32'
Set-TextValue $ws.Range("E23") '  +3.01%  '

Set-TextValue $ws.Range("D24") '104.66'
Set-TextValue $ws.Range("E24") '  +2.98%  '

Set-TextValue $ws.Range("E25") '  -5.35%  '

Set-TextValue $ws.Range("D26") '2.89'
Set-TextValue $ws.Range("E26") '  -4.31%  '

Set-TextValue $ws.Range("D27") '10.15'
Set-TextValue $ws.Range("E27") '  -5.63%  '

Set-TextValue $ws.Range("D28") '9.59'
Set-TextValue $ws.Range("E28") '  -0.03%  '

Set-TextValue $ws.Range("D29") '33.13'
Set-TextValue $ws.Range("E29") '  -2.33%  '

Set-TextValue $ws.Range("D30") '4.39'
Set-TextValue $ws.Range("E30") '  -7.77%  '

Set-TextValue $ws.Range("D31") '6.76'
Set-TextValue $ws.Range("E31") '  -5.01%  '

Set-TextValue $ws.Range("D32") '11.74'
Set-TextValue $ws.Range("E32") '  -4.57%  '

Set-TextValue $ws.Range("E33") '  -5.34%  '

Set-TextValue $ws.Range("D34") '61.94'
Set-TextValue $ws.Range("E34") '  -2.04%  '

Set-TextValue $ws.Range("D35") '3.789.70'
Set-TextValue $ws.Range("E35") '  -4.23%  '

$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D36") '0.0₃0812'
Set-TextValue $ws.Range("E36") '  -9.28%  '

$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D37") '3.72'
Set-TextValue $ws.Range("E37") '  +5.22%  '

Set-TextValue $ws.Range("D38") '1.00'
Set-TextValue $ws.Range("E38") '  +0.02%  '

Set-TextValue $ws.Range("D39") '501.45'
Set-TextValue $ws.Range("E39") '  -4.09%  '

Set-TextValue $ws.Range("D40") '2.91'
Set-TextValue $ws.Range("E40") '  -6.57%  '

Set-TextValue $ws.Range("D41") '0.371'
Set-TextValue $ws.Range("E41") '  -5.08%  '

Set-TextValue $ws.Range("E42") '  +0.32%  '

Set-TextValue $ws.Range("D43") '34.44'
Set-TextValue $ws.Range("E43") '  -6.61%  '

Set-TextValue $ws.Range("D44") '0.0447'
Set-TextValue $ws.Range("E44") '  -1.72%  '

Set-TextValue $ws.Range("E45") '  -0.07%  '

Set-TextValue $ws.Range("E46") '  -0.35%  '

Set-TextValue $ws.Range("E47") '  -3.18%  '

Set-TextValue $ws.Range("E48") '  -0.10%  '

Set-TextValue $ws.Range("E49") '  -3.83%  '

Set-TextValue $ws.Range("D50") '1.76'
Set-TextValue $ws.Range("E50") '  +19.14%  '

$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue $ws.Range("D51") '0.000233'
Set-TextValue $ws.Range("E51") '  -7.48%  '
